# Updates the options-Greeks columns (P..T) on the active sheet.
# For every data row (2-37):
#   - Gamma (Q) and Theta (T) are cleared to an empty text value
#     (matches the exported "blank" cell produced by the upstream
#     pandas/openpyxl pipeline - an explicit empty string, not a
#     fully-deleted cell).
#   - Rho (S) is reset to 0 for the rows that still carried a stale
#     non-zero value (rows 2-17).
#   - Delta (P) and Vega (R) are reset to 0 for the two rows that still
#     carried stale non-zero values (rows 18 and 22).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Clear-ToEmptyText($cell) {
    # Leading apostrophe forces a genuine (empty) text cell instead of
    # Excel's usual "blank -> delete the cell" behaviour; re-applying the
    # Normal style afterwards drops the quote-prefix formatting flag that
    # the apostrophe trick leaves behind so no visible formatting changes.
    $cell.Formula = "'"
    $cell.Style = "Normal"
}

for ($row = 2; $row -le 37; $row++) {
    Clear-ToEmptyText $ws.Range("Q$row")
    Clear-ToEmptyText $ws.Range("T$row")

    if ($row -le 17) {
        $ws.Range("S$row").Value = 0
    }

    if ($row -eq 18 -or $row -eq 22) {
        $ws.Range("P$row").Value = 0
        $ws.Range("R$row").Value = 0
    }
}
